$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated B:G values for rows 2-6 (regen sval data to filter save games)
$ws.Range("B2").Value = 0.3048080303191223
$ws.Range("C2").Value = 0.3127903958511391
$ws.Range("D2").Value = 0.8054896365839992
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 1.919867272924993

$ws.Range("B3").Value = 0.3048080303191223
$ws.Range("C3").Value = 1.667794583268128
$ws.Range("D3").Value = 0.8054896365839992
$ws.Range("E3").Value = 8.660232485948974
$ws.Range("G3").Value = 11.43832473612022

$ws.Range("B4").Value = 0.3048080303191223
$ws.Range("C4").Value = 114.8270160096505
$ws.Range("D4").Value = 26.21740644021617
$ws.Range("E4").Value = 645.3272768299601
$ws.Range("G4").Value = 786.6765073101459

$ws.Range("B5").Value = 0.3048080303191223
$ws.Range("C5").Value = 0.3127903958511391
$ws.Range("D5").Value = 3.900430680208489
$ws.Range("E5").Value = 616238.5361209477
$ws.Range("G5").Value = 616243.0541500541

$ws.Range("B6").Value = 1.459612070389937
$ws.Range("C6").Value = 1.667794583268128
$ws.Range("D6").Value = 0.1575252929769615
$ws.Range("E6").Value = 8.660232485948974
$ws.Range("G6").Value = 11.945164432584
